# Automatische test-sync: 2025-08-19 21:03:50
# Adds the new mail-log entry (row 22) on the "Logs" sheet, extends the
# conditional-formatting ranges to include the new row, and refreshes the
# aggregate count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append the new row ---------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A22").Value = "Vraag over product"
$logs.Range("B22").Value = "documentatie@testbedrijf123.nl"
$logs.Range("D22").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("F22").Value = "2025-08-19 21:02:55"
$logs.Range("G22").Value = "Nee"
$logs.Range("H22").Value = "Ja"
$logs.Range("I22").Value = "Nee"
$logs.Range("J22").Value = "Nee"

# --- Extend conditional formatting ranges to cover the new row --------
$logs.Range("D2:D21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D22"))
$logs.Range("G2:G21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G22"))
$logs.Range("H2:H21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H22"))
$logs.Range("I2:I21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I22"))
$logs.Range("J2:J21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J22"))

# --- "Dashboard" sheet: bump the aggregate count -----------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 21
